# Update the "Share Value" column (B) with newly regenerated secret shares,
# and extend the secret-shares table with 10 additional participants
# (Person Id 31-40), matching the refreshed share data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shares = @(
    @(1, 27446547),
    @(2, 11284996),
    @(3, 27237078),
    @(4, 7025799),
    @(5, 12207958),
    @(6, 32847345),
    @(7, 21069208),
    @(8, 25278121),
    @(9, 22385649),
    @(10, 18473749),
    @(11, 15240303),
    @(12, 9999118),
    @(13, 29234393),
    @(14, 27937318),
    @(15, 23823942),
    @(16, 30226239),
    @(17, 22537641),
    @(18, 5321972),
    @(19, 12313448),
    @(20, 5753276),
    @(21, 10607522),
    @(22, 13903710),
    @(23, 31839756),
    @(24, 9120567),
    @(25, 20284843),
    @(26, 1269341),
    @(27, 17844611),
    @(28, 30733727),
    @(29, 29830155),
    @(30, 643286),
    @(31, 24961837),
    @(32, 12418115),
    @(33, 2478220),
    @(34, 30224177),
    @(35, 25690535),
    @(36, 15191169),
    @(37, 20655879),
    @(38, 26075923),
    @(39, 11058484),
    @(40, 17935604)
)

foreach ($row in $shares) {
    $personId = $row[0]
    $shareValue = $row[1]
    $excelRow = $personId + 1

    $ws.Cells.Item($excelRow, 1).Value = $personId
    $ws.Cells.Item($excelRow, 2).Value = $shareValue
}
